# Region II_GABALDON.xlsx edit:
# Insert a new "INDEX (DO NOT MODIFY)" column at the very left of the data
# sheet, renumbering/shifting every existing column one place to the right,
# fill in the new index values per row, and re-case several header labels
# to ALL CAPS (matching the republished template).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Insert new column A, pushing A:AB -> B:AC -------------------------
$ws.Columns("A:A").Insert()

# Match the new column's width to the OOXML width="23" used by the template
# (COM ColumnWidth is offset from the OOXML character-width units by the
# sheet's default-font padding, ~0.83 for Calibri 11 in this workbook).
$ws.Columns("A:A").ColumnWidth = 22.17

# --- 2. Clone formatting into the new column from its neighbours ----------
# Header cell (row 1) takes the header style; body cells take the body style.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. New column header + per-row index values ---------------------------
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"

$ws.Range("A2").Value = 21
$ws.Range("A3").Value = 22
$ws.Range("A4").Value = 23
$ws.Range("A5").Value = 24
$ws.Range("A6").Value = 139
$ws.Range("A7").Value = 140
$ws.Range("A8").Value = 141
$ws.Range("A9").Value = 172
$ws.Range("A10").Value = 187
$ws.Range("A11").Value = 223

# --- 4. Re-cased / re-worded header labels (now shifted one column right) -
$ws.Range("C1").Value = "REGION"
$ws.Range("R1").Value = " TARGET COMPLETION DATE "
$ws.Range("S1").Value = "ACTUAL DATE OF COMPLETION"
$ws.Range("T1").Value = "PROJECT ID"
$ws.Range("U1").Value = "CONTRACT ID"
$ws.Range("V1").Value = "ISSUANCE OF INVITATION TO BID"
$ws.Range("W1").Value = "PRE-SUBMISSION CONFERENCE"
$ws.Range("X1").Value = "BID OPENING"
$ws.Range("Y1").Value = "ISSUANCE OF RESOLUTION TO AWARD"
$ws.Range("Z1").Value = "ISSUANCE OF NOTICE TO PROCEED"
$ws.Range("AA1").Value = "NAME OF CONTRACTOR"
$ws.Range("AB1").Value = "OTHER REMARKS"
